$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 429.21054
$ws.Range("I28").Value = 437
$ws.Range("K28").Value = 437
$ws.Range("M28").Value = 48
$ws.Range("H32").Value = 31251768
$ws.Range("J32").Value = 35716108
$ws.Range("L32").Value = 35716108
$ws.Range("N32").Value = -35716760
$ws.Range("H40").Value = 2728.3333
$ws.Range("J40").Value = 2314.5
$ws.Range("L40").Value = 2314.5
$ws.Range("N40").Value = -2664.5
$ws.Range("H62").Value = 4845
$ws.Range("I62").Value = 4837.143
$ws.Range("K62").Value = 4837.143
$ws.Range("M62").Value = -4213.143
$ws.Range("H65").Value = 4845
$ws.Range("I65").Value = 4837.143
$ws.Range("K65").Value = 24185.715
$ws.Range("M65").Value = -21065.715
$ws.Range("H70").Value = 1129.5
$ws.Range("J70").Value = 1262
$ws.Range("L70").Value = 3786
$ws.Range("N70").Value = -4326
$ws.Range("H73").Value = 1129.5
$ws.Range("J73").Value = 1262
$ws.Range("L73").Value = 3786
$ws.Range("N73").Value = -5658
$ws.Range("H76").Value = 3709.1667
$ws.Range("I76").Value = 3482
$ws.Range("J76").Value = 4163.5
$ws.Range("K76").Value = 3482
$ws.Range("L76").Value = 4163.5
$ws.Range("M76").Value = -3167
$ws.Range("N76").Value = -4793.5
$ws.Range("H79").Value = 3709.1667
$ws.Range("I79").Value = 3482
$ws.Range("J79").Value = 4163.5
$ws.Range("K79").Value = 3482
$ws.Range("L79").Value = 4163.5
$ws.Range("M79").Value = -2390
$ws.Range("N79").Value = -6347.5
$ws.Range("H88").Value = 6580.8
$ws.Range("I88").Value = 798
$ws.Range("J88").Value = 15255
$ws.Range("K88").Value = 798
$ws.Range("L88").Value = 15255
$ws.Range("M88").Value = -392
$ws.Range("N88").Value = -16067
$ws.Range("H91").Value = 6580.8
$ws.Range("I91").Value = 798
$ws.Range("J91").Value = 15255
$ws.Range("K91").Value = 798
$ws.Range("L91").Value = 15255
$ws.Range("M91").Value = 606
$ws.Range("N91").Value = -18063
$ws.Range("H92").Value = 3348898
$ws.Range("I92").Value = 1563227.8
$ws.Range("J92").Value = 7813073.5
$ws.Range("K92").Value = 1563227.8
$ws.Range("L92").Value = 7813073.5
$ws.Range("M92").Value = -1561979.8
$ws.Range("N92").Value = -7815569.5
$ws.Range("H107").Value = 13890434
$ws.Range("I107").Value = 1209.6471
$ws.Range("K107").Value = 1209.6471
$ws.Range("M107").Value = 710.3529000000001
$ws.Range("H112").Value = 72430.87
$ws.Range("J112").Value = 72430.87
$ws.Range("L112").Value = 217292.61
$ws.Range("N112").Value = -219508.61
$ws.Range("H113").Value = 3751.5334
$ws.Range("I113").Value = 3826.7778
$ws.Range("K113").Value = 3826.7778
$ws.Range("M113").Value = -572.7777999999998
$ws.Range("H115").Value = 1999
$ws.Range("I115").Value = 1999
$ws.Range("K115").Value = 5997
$ws.Range("M115").Value = -4430
$ws.Range("H116").Value = 6282.75
$ws.Range("J116").Value = 6666
$ws.Range("L116").Value = 6666
$ws.Range("N116").Value = -13550
$ws.Range("H137").Value = 2711233
$ws.Range("I137").Value = 4831.769
$ws.Range("K137").Value = 14495.307
$ws.Range("M137").Value = -11945.307
$ws.Range("H138").Value = 5382.711
$ws.Range("J138").Value = 2440.577
$ws.Range("L138").Value = 7321.731000000001
$ws.Range("N138").Value = -17601.731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3216.4102
$ws.Range("I32").Value = 547.86664
$ws.Range("J32").Value = 12111.556
$ws.Range("K32").Value = 547.86664
$ws.Range("L32").Value = 12111.556
$ws.Range("M32").Value = -260.86664
$ws.Range("N32").Value = -12685.556
$ws.Range("H61").Value = 1526577.6
$ws.Range("I61").Value = 49204.176
$ws.Range("J61").Value = 3525377
$ws.Range("K61").Value = 49204.176
$ws.Range("L61").Value = 3525377
$ws.Range("M61").Value = -48992.176
$ws.Range("N61").Value = -3525801
$ws.Range("H74").Value = 402809.47
$ws.Range("I74").Value = 2154.7896
$ws.Range("J74").Value = 1164053.4
$ws.Range("K74").Value = 2154.7896
$ws.Range("L74").Value = 1164053.4
$ws.Range("M74").Value = -1280.7896
$ws.Range("N74").Value = -1165801.4
$ws.Range("H77").Value = 402809.47
$ws.Range("I77").Value = 2154.7896
$ws.Range("J77").Value = 1164053.4
$ws.Range("K77").Value = 10773.948
$ws.Range("L77").Value = 5820267
$ws.Range("M77").Value = -6405.948
$ws.Range("N77").Value = -5829003
$ws.Range("H88").Value = 2490.9167
$ws.Range("J88").Value = 2499.182
$ws.Range("L88").Value = 2499.182
$ws.Range("N88").Value = -3311.182
$ws.Range("H91").Value = 2490.9167
$ws.Range("J91").Value = 2499.182
$ws.Range("L91").Value = 2499.182
$ws.Range("N91").Value = -5307.182
$ws.Range("H97").Value = 9449.538
$ws.Range("I97").Value = 10415.818
$ws.Range("J97").Value = 4135
$ws.Range("K97").Value = 10415.818
$ws.Range("L97").Value = 4135
$ws.Range("M97").Value = -9919.817999999999
$ws.Range("N97").Value = -5127
$ws.Range("H109").Value = 89688.5
$ws.Range("J109").Value = 89688.5
$ws.Range("L109").Value = 89688.5
$ws.Range("N109").Value = -92462.5
$ws.Range("H122").Value = 1557.4445
$ws.Range("J122").Value = 2468
$ws.Range("L122").Value = 7404
$ws.Range("N122").Value = -12304
$ws.Range("H132").Value = 2831.5264
$ws.Range("I132").Value = 1899.8572
$ws.Range("J132").Value = 3982.4119
$ws.Range("K132").Value = 5699.571599999999
$ws.Range("L132").Value = 11947.2357
$ws.Range("M132").Value = -3169.571599999999
$ws.Range("N132").Value = -17007.2357
$ws.Range("H133").Value = 103163
$ws.Range("J133").Value = 103163
$ws.Range("L133").Value = 103163
$ws.Range("N133").Value = -108223
$ws.Range("H136").Value = 1526577.6
$ws.Range("I136").Value = 49204.176
$ws.Range("J136").Value = 3525377
$ws.Range("K136").Value = 147612.528
$ws.Range("L136").Value = 10576131
$ws.Range("M136").Value = -145062.528
$ws.Range("N136").Value = -10581231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H59").Value = 110000
$ws.Range("I59").Value = 100000
$ws.Range("K59").Value = 100000
$ws.Range("M59").Value = -99153
$ws.Range("H99").Value = 20653.215
$ws.Range("I99").Value = 23154.223
$ws.Range("K99").Value = 23154.223
$ws.Range("M99").Value = -21656.223
$ws.Range("H105").Value = 18868.75
$ws.Range("I105").Value = 15142.6
$ws.Range("K105").Value = 15142.6
$ws.Range("M105").Value = -13395.6
$ws.Range("H133").Value = 58330
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 58330
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 58330
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -68450
$ws.Range("H134").Value = 37502116
$ws.Range("I134").Value = 2442.1667
$ws.Range("K134").Value = 7326.500100000001
$ws.Range("M134").Value = -4791.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000500
$ws.Range("H31").Value = 2726.6099
$ws.Range("I31").Value = 2943
$ws.Range("K31").Value = 2943
$ws.Range("M31").Value = -2648
$ws.Range("H34").Value = 2726.6099
$ws.Range("I34").Value = 2943
$ws.Range("K34").Value = 2943
$ws.Range("M34").Value = -2741
$ws.Range("H93").Value = 15995.223
$ws.Range("I93").Value = 10494.625
$ws.Range("K93").Value = 10494.625
$ws.Range("M93").Value = -8622.625
$ws.Range("H122").Value = 2200.6667
$ws.Range("I122").Value = 2101.5715
$ws.Range("J122").Value = 2547.5
$ws.Range("K122").Value = 6304.7145
$ws.Range("L122").Value = 7642.5
$ws.Range("M122").Value = -3854.7145
$ws.Range("N122").Value = -12542.5
$ws.Range("H133").Value = 83119.664
$ws.Range("J133").Value = 83119.664
$ws.Range("L133").Value = 83119.664
$ws.Range("N133").Value = -88179.664
$ws.Range("H134").Value = 2218.9375
$ws.Range("I134").Value = 2069.5454
$ws.Range("J134").Value = 2547.6
$ws.Range("K134").Value = 6208.6362
$ws.Range("L134").Value = 7642.799999999999
$ws.Range("M134").Value = -3673.6362
$ws.Range("N134").Value = -12712.8
$ws.Range("H141").Value = 210624
$ws.Range("J141").Value = 210624
$ws.Range("L141").Value = 210624
$ws.Range("N141").Value = -220984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1072.1842
$ws.Range("J2").Value = 1462.4348
$ws.Range("L2").Value = 8774.6088
$ws.Range("N2").Value = -9000.6088
$ws.Range("H4").Value = 877495.9
$ws.Range("I4").Value = 797771.9
$ws.Range("K4").Value = 2393315.7
$ws.Range("M4").Value = -2393203.7
$ws.Range("H25").Value = 1974.75
$ws.Range("J25").Value = 2950
$ws.Range("L25").Value = 8850
$ws.Range("N25").Value = -9188
$ws.Range("H30").Value = 1974.75
$ws.Range("J30").Value = 2950
$ws.Range("L30").Value = 8850
$ws.Range("N30").Value = -9054
$ws.Range("H39").Value = 8352.429
$ws.Range("J39").Value = 11233.6
$ws.Range("L39").Value = 33700.8
$ws.Range("N39").Value = -34288.8
$ws.Range("H47").Value = 1642.8572
$ws.Range("I47").Value = 1000
$ws.Range("J47").Value = 1692.3077
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 5076.9231
$ws.Range("M47").Value = -2569
$ws.Range("N47").Value = -5938.9231
$ws.Range("H57").Value = 8180
$ws.Range("I57").Value = 4825
$ws.Range("J57").Value = 10416.667
$ws.Range("K57").Value = 14475
$ws.Range("L57").Value = 31250.001
$ws.Range("M57").Value = -13916
$ws.Range("N57").Value = -32368.001
$ws.Range("H80").Value = 2166.3333
$ws.Range("I80").Value = 2249.5
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 6748.5
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -5812.5
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2166.3333
$ws.Range("I83").Value = 2249.5
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 20245.5
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -15565.5
$ws.Range("N83").Value = -27360
$ws.Range("H92").Value = 446.75
$ws.Range("I92").Value = 465.66666
$ws.Range("J92").Value = 390
$ws.Range("K92").Value = 1396.99998
$ws.Range("L92").Value = 1170
$ws.Range("M92").Value = -148.9999800000001
$ws.Range("N92").Value = -3666
$ws.Range("H117").Value = 74076830
$ws.Range("J117").Value = 74076830
$ws.Range("L117").Value = 222230490
$ws.Range("N117").Value = -222237374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 126.875
$ws.Range("I2").Value = 113.833336
$ws.Range("J2").Value = 166
$ws.Range("K2").Value = 113.833336
$ws.Range("L2").Value = 166
$ws.Range("M2").Value = -0.8333360000000027
$ws.Range("N2").Value = -392
$ws.Range("H15").Value = 35593.668
$ws.Range("J15").Value = 35593.668
$ws.Range("L15").Value = 35593.668
$ws.Range("N15").Value = -36169.668
$ws.Range("H35").Value = 24980
$ws.Range("I35").Value = 24980
$ws.Range("K35").Value = 24980
$ws.Range("M35").Value = -24682
$ws.Range("H42").Value = 71999.664
$ws.Range("J42").Value = 71999.664
$ws.Range("L42").Value = 71999.664
$ws.Range("N42").Value = -72969.664
$ws.Range("H43").Value = 2824.5
$ws.Range("I43").Value = 2824.5
$ws.Range("K43").Value = 2824.5
$ws.Range("M43").Value = -2673.5
$ws.Range("H62").Value = 47548.43
$ws.Range("I62").Value = 41947
$ws.Range("K62").Value = 41947
$ws.Range("M62").Value = -41261
$ws.Range("H65").Value = 47548.43
$ws.Range("I65").Value = 41947
$ws.Range("K65").Value = 125841
$ws.Range("M65").Value = -122409
$ws.Range("H68").Value = 50000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H70").Value = 5127
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5127
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5127
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5667
$ws.Range("H71").Value = 50000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H73").Value = 5127
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5127
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5127
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6999
$ws.Range("H81").Value = 35593.668
$ws.Range("J81").Value = 35593.668
$ws.Range("L81").Value = 35593.668
$ws.Range("N81").Value = -37589.668
$ws.Range("H84").Value = 35593.668
$ws.Range("J84").Value = 35593.668
$ws.Range("L84").Value = 106781.004
$ws.Range("N84").Value = -116765.004
$ws.Range("H97").Value = 183017.45
$ws.Range("I97").Value = 333929
$ws.Range("J97").Value = 126425.625
$ws.Range("K97").Value = 333929
$ws.Range("L97").Value = 126425.625
$ws.Range("M97").Value = -333433
$ws.Range("N97").Value = -127417.625
$ws.Range("H106").Value = 36244.75
$ws.Range("J106").Value = 36244.75
$ws.Range("L106").Value = 36244.75
$ws.Range("N106").Value = -38768.75
$ws.Range("H113").Value = 9767.143
$ws.Range("I113").Value = 9897
$ws.Range("J113").Value = 9594
$ws.Range("K113").Value = 9897
$ws.Range("L113").Value = 9594
$ws.Range("M113").Value = -7727
$ws.Range("N113").Value = -13934
$ws.Range("H115").Value = 71999.664
$ws.Range("J115").Value = 71999.664
$ws.Range("L115").Value = 71999.664
$ws.Range("N115").Value = -74349.664
$ws.Range("H137").Value = 184584
$ws.Range("J137").Value = 189479.33
$ws.Range("L137").Value = 189479.33
$ws.Range("N137").Value = -199679.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 979.7619
$ws.Range("I16").Value = 865.2778
$ws.Range("K16").Value = 865.2778
$ws.Range("M16").Value = -695.2778
$ws.Range("H22").Value = 4509.409
$ws.Range("I22").Value = 1320.5
$ws.Range("J22").Value = 5705.25
$ws.Range("K22").Value = 1320.5
$ws.Range("L22").Value = 5705.25
$ws.Range("M22").Value = -1025.5
$ws.Range("N22").Value = -6295.25
$ws.Range("H27").Value = 4509.409
$ws.Range("I27").Value = 1320.5
$ws.Range("J27").Value = 5705.25
$ws.Range("K27").Value = 1320.5
$ws.Range("L27").Value = 5705.25
$ws.Range("M27").Value = -1213.5
$ws.Range("N27").Value = -5919.25
$ws.Range("H35").Value = 6780
$ws.Range("I35").Value = 4475
$ws.Range("J35").Value = 16000
$ws.Range("K35").Value = 4475
$ws.Range("L35").Value = 16000
$ws.Range("M35").Value = -4139
$ws.Range("N35").Value = -16672
$ws.Range("H61").Value = 2888.2104
$ws.Range("I61").Value = 2062.6428
$ws.Range("J61").Value = 5199.8
$ws.Range("K61").Value = 2062.6428
$ws.Range("L61").Value = 5199.8
$ws.Range("M61").Value = -1860.6428
$ws.Range("N61").Value = -5603.8
$ws.Range("H93").Value = 4227.857
$ws.Range("I93").Value = 2432.6667
$ws.Range("J93").Value = 14999
$ws.Range("K93").Value = 2432.6667
$ws.Range("L93").Value = 14999
$ws.Range("M93").Value = -1184.6667
$ws.Range("N93").Value = -17495
$ws.Range("H100").Value = 4305.6875
$ws.Range("I100").Value = 3878
$ws.Range("J100").Value = 4638.3335
$ws.Range("K100").Value = 3878
$ws.Range("L100").Value = 4638.3335
$ws.Range("M100").Value = -3337
$ws.Range("N100").Value = -5720.3335
$ws.Range("H110").Value = 41713.145
$ws.Range("J110").Value = 41713.145
$ws.Range("L110").Value = 41713.145
$ws.Range("N110").Value = -49893.145
$ws.Range("H113").Value = 2888.2104
$ws.Range("I113").Value = 2062.6428
$ws.Range("J113").Value = 5199.8
$ws.Range("K113").Value = 2062.6428
$ws.Range("L113").Value = 5199.8
$ws.Range("M113").Value = 107.3571999999999
$ws.Range("N113").Value = -9539.799999999999
$ws.Range("H132").Value = 4158.5
$ws.Range("J132").Value = 5498.3335
$ws.Range("L132").Value = 16495.0005
$ws.Range("N132").Value = -21555.0005
$ws.Range("H135").Value = 163495
$ws.Range("J135").Value = 163495
$ws.Range("L135").Value = 163495
$ws.Range("N135").Value = -173635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 75000
$ws.Range("J16").Value = 75000
$ws.Range("L16").Value = 75000
$ws.Range("N16").Value = -75584
$ws.Range("H62").Value = 2433.4
$ws.Range("J62").Value = 2333.4
$ws.Range("L62").Value = 2333.4
$ws.Range("N62").Value = -3581.4
$ws.Range("H65").Value = 2433.4
$ws.Range("J65").Value = 2333.4
$ws.Range("L65").Value = 11667
$ws.Range("N65").Value = -17907
$ws.Range("H96").Value = 27899.5
$ws.Range("I96").Value = 3866
$ws.Range("K96").Value = 3866
$ws.Range("M96").Value = -2493
$ws.Range("H107").Value = 2199669.8
$ws.Range("I107").Value = 1422.125
$ws.Range("K107").Value = 4266.375
$ws.Range("M107").Value = -2346.375
$ws.Range("H113").Value = 478.80646
$ws.Range("I113").Value = 419.625
$ws.Range("K113").Value = 1258.875
$ws.Range("M113").Value = 911.125
$ws.Range("H122").Value = 1590.5385
$ws.Range("I122").Value = 1590.5385
$ws.Range("K122").Value = 4771.6155
$ws.Range("M122").Value = -2321.6155
$ws.Range("H126").Value = 20835714
$ws.Range("I126").Value = 50002296
$ws.Range("K126").Value = 150006888
$ws.Range("M126").Value = -150004418
$ws.Range("H131").Value = 39143
$ws.Range("J131").Value = 39143
$ws.Range("L131").Value = 39143
$ws.Range("N131").Value = -49223
